# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    3  = 2959
    4  = 191
    7  = 1613
    10 = 26
    11 = 1321
    13 = 438
    15 = 68
    16 = 57
    17 = 118
    19 = 99
    20 = 3032
    21 = 369
    22 = 79
    23 = 19
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
